# Actualización automática 2025-09-26 14:45:08
$wb = $excel.ActiveWorkbook

# Sheet: VENTAS POR GRUPO
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M13").Value = 11932.44

# Sheet: VENTA MENSUAL
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F13").Value = 13291.67
$ws2.Range("F23").Value = 52766.17

# Sheet: CUMPLIMIENTO MENSUAL
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D12").Value = 45557.33
$ws3.Range("E12").Value = -8733.686907882904
$ws3.Range("F12").Value = 1.237176068810871

$ws3.Range("D15").Value = 52323.03
$ws3.Range("E15").Value = 3101.713166133772
$ws3.Range("F15").Value = 0.9440373921655082
